$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")
$ws.Activate()

# Scrum daily meeting update: status column (F) changes for two tasks.
# Row 7 "Review code" moves from "In progress" to "Done".
$ws.Range("F7").Value = "Done"
# Row 13 "Document implementation" moves from "To do" to "In progress".
$ws.Range("F13").Value = "In progress"

# Day 11 / Day 12 daily-meeting effort logged (column R = Day 12) for the
# tasks that had progress reported today.
$ws.Range("R7").Value = 1
$ws.Range("R8").Value = 1
$ws.Range("R11").Value = 1
$ws.Range("R13").Value = 1

# Update the active selection to reflect where the user left off.
$ws.Range("R11").Select()
